$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rework the tail of the "proposed solution" paragraph:
#    - "...every order. A" + " feature..." -> "...every order. Additionally, a feature..."
#    - "...part of the Sun Room with some renovations into such a facility."
#        -> "...part of the Sun Room and some renovations to transform it into
#            such a facility. " (note trailing space)
# ---------------------------------------------------------------------------
$oldTail = "order. A feature that UDCC has is the special dietary kitchen for those people who have extremely inhibiting allergies and allows for a sterile place for them to get and eat food. Seasons should also adopt this idea, and this could be easily achieved by reallocating part of the Sun Room with some renovations into such a facility."
$newTail = "order. Additionally, a feature that UDCC has is the special dietary kitchen for those people who have extremely inhibiting allergies and allows for a sterile place for them to get and eat food. Seasons should also adopt this idea, and this could be easily achieved by reallocating part of the Sun Room and some renovations to transform it into such a facility. "

$null = $d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)

# ---------------------------------------------------------------------------
# 2. Re-establish run boundaries inside that rewritten sentence so the XML
#    ends up with the same granular run layout as the authored edit (rather
#    than one big merged run). Each segment below becomes its own <w:r> by
#    toggling a character property on/off across its exact span.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$anchor = $full.IndexOf("order. Additionally, a feature")
$start = $anchor + "order. ".Length

$segments = @(
    "A",
    "dditionally, a",
    " feature that UDCC has is the special dietary kitchen for those people who have ",
    "extremely inhibiting allergies and allows for a sterile place for them to get and eat food. Seasons should also adopt this idea, and this could be easily achieved by reallocating part of the ",
    "Sun Room",
    " ",
    "and",
    " some renovations ",
    "to transform it ",
    "into such a facility",
    ".",
    " "
)

$pos = $start
foreach ($seg in $segments) {
    $segStart = $pos
    $segEnd = $pos + $seg.Length
    $segRange = $d.Range($segStart, $segEnd)
    $segRange.Font.Bold = $true
    $segRange.Font.Bold = $false
    $pos = $segEnd
}

# ---------------------------------------------------------------------------
# 3. Append the new material at the end of the document:
#      <blank paragraph>
#      2.2 Result        (italic, underlined - matches the other "x.y Title"
#                          headings already used throughout the proposal)
#      <tab>The
#
#    All the new paragraphs/text are created first (with default / inherited
#    formatting), and only once every paragraph exists do we go back and
#    apply the italic+underline styling to the "2.2 Result" heading. Doing
#    the styling while the heading paragraph is still the very last
#    paragraph in the document causes the engine to "leak" that formatting
#    into whatever paragraph is created next, so it has to happen last.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Last
$blankRange = $blankPara.Range
$blankRange.Collapse(0)
$blankRange.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Last
$headingRange = $headingPara.Range
$headingRange.Collapse(0)
$headingRange.InsertAfter("2.2 Result")

$headingEndRange = $headingPara.Range
$headingEndRange.Collapse(0)
$headingEndRange.InsertParagraphAfter()

$finalPara = $d.Paragraphs.Last
$finalRange = $finalPara.Range
$finalRange.Collapse(0)
$finalRange.InsertAfter([char]9 + "The ")

# -- now, with every paragraph already in place, style the heading text --
$full2 = $d.Content.Text
$headingIdx = $full2.IndexOf("2.2 Result")
$headingStart = $headingIdx
$headingEnd = $headingIdx + "2.2 Result".Length

$headingTextRange = $d.Range($headingStart, $headingEnd)
$headingTextRange.Font.Italic = $true
$headingTextRange.Font.ItalicBi = $true
$headingTextRange.Font.Underline = 1

$headingParaMark = $headingPara.Range
$headingParaMark.Collapse(0)
$headingParaMark.Font.Italic = $true
$headingParaMark.Font.ItalicBi = $true
$headingParaMark.Font.Underline = 1

Write-Output ("paragraph count = " + $d.Paragraphs.Count)
